# Edit the mTBI conditions sheet:
#  - Remove the "opacity" column (D) entirely
#  - Change column C values to a fixed rating question text
#  - Bold the new column C values (rows 2-7)
#  - Update image filename values to include an "images/" prefix

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear column D (removes the "opacity" header + the "1" values) and shrink
# the used range down to A1:C7.
$ws.Range("D1:D7").Clear()

# Header row stays the same for A:C, just drop D.
$ws.Range("A1").Value = "arrow"
$ws.Range("B1").Value = "image"
$ws.Range("C1").Value = "q_or_r"

$question = "Rate your ability to control your brain"

$ws.Range("A2").Value = "down"
$ws.Range("B2").Value = "images/down.jpg"
$ws.Range("C2").Value = $question

$ws.Range("A3").Value = "up"
$ws.Range("B3").Value = "images/up.jpg"
$ws.Range("C3").Value = $question

$ws.Range("A4").Value = "up"
$ws.Range("B4").Value = "images/up.jpg"
$ws.Range("C4").Value = $question

$ws.Range("A5").Value = "down"
$ws.Range("B5").Value = "images/down.jpg"
$ws.Range("C5").Value = $question

$ws.Range("A6").Value = "down"
$ws.Range("B6").Value = "images/down.jpg"
$ws.Range("C6").Value = $question

$ws.Range("A7").Value = "up"
$ws.Range("B7").Value = "images/up.jpg"
$ws.Range("C7").Value = $question

# Bold the rating-question column for the data rows.
$ws.Range("C2:C7").Font.Bold = $true

# Selection ends up on C11 in the final file.
$ws.Range("C11").Select() | Out-Null

# Column widths for the data sheet (A/B/C get explicit widths, rest default).
$ws.Columns.Item(1).ColumnWidth = 11.9843137254902
$ws.Columns.Item(2).ColumnWidth = 15.121568627451
$ws.Columns.Item(3).ColumnWidth = 31.321568627451

# The other two (empty) sheets also pick up the refreshed default column width.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Columns.Item(1).ColumnWidth = 11.9843137254902

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Columns.Item(1).ColumnWidth = 11.9843137254902
